# Apply the "add monster and cavalry classes" edit to battle_5.xlsx
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "nodes": two group_3 cells get populated with tunnel refs
# ---------------------------------------------------------------
$nodes = $wb.Worksheets.Item("nodes")
$nodes.Range("F32").Copy($nodes.Range("G32"))
$nodes.Range("G32").Value = "tunnel_9"

$nodes.Range("E82").Copy($nodes.Range("G82"))
$nodes.Range("G82").Value = "tunnel_2"

# ---------------------------------------------------------------
# Sheet "parameters": add the two new cavalry parameters, right
# before the flier parameters
# ---------------------------------------------------------------
$params = $wb.Worksheets.Item("parameters")
$params.Range("A8:B9").EntireRow.Insert()
$params.Range("A8").Value = "cavalry_distance"
$params.Range("B8").Value = 4.5
$params.Range("A9").Value = "cavalry_height_difference_threshold"
$params.Range("B9").Value = 2

# ---------------------------------------------------------------
# Sheet "interactions": add a "cavalry" interaction column (copy of
# the existing siege behaviour), collapse the per-tunnel rows into a
# single "all tunnels" row.
# ---------------------------------------------------------------
$inter = $wb.Worksheets.Item("interactions")

# insert a new column before the old "siege" column (F) and fill it
# in with the same values siege has (cavalry behaves like siege)
$inter.Columns.Item(6).Insert()
$inter.Range("F1").Value = "cavalry"
for ($r = 2; $r -le 21; $r++) {
    $v = $inter.Cells.Item($r, 7).Value2
    $inter.Cells.Item($r, 6).Value = $v
}

# collapse the individual tunnel_x / tunnel_crossing rows (13-21) into
# row 12, which now represents the full tunnel list
$inter.Range("A13:G21").EntireRow.Delete()
$inter.Range("A12").Value = "tunnel_1,tunnel_2,tunnel_3,tunnel_4,tunnel_5,tunnel_6,tunnel_7,tunnel_8,tunnel_9,tunnel_crossing"
$inter.Range("B12").Value = "tunnel_1,tunnel_2,tunnel_3,tunnel_4,tunnel_5,tunnel_6,tunnel_7,tunnel_8,tunnel_9,tunnel_crossing"

# ---------------------------------------------------------------
# Restore selection / active-sheet state to match the edited file
# ---------------------------------------------------------------
$units = $wb.Worksheets.Item("units")
$units.Range("F16").Select()

$params.Range("B13").Select()

$nodes.Activate()
$nodes.Range("L35").Select()

$inter.Activate()
$inter.Range("D30").Select()
